$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "Или когда настроение" "Или в плохом настроении"
Replace-Text "Честь растирая" "И честь растирая"
Replace-Text "Не тем, как деды" "Не прошлым - как деды"
Replace-Text "А тем, как братья мои" "А тем, что братья мои"
Replace-Text "И просто, без всяких затей" "И просто - без всяких затей"
